# November 2014 payslips - add department-driven HRA and Food Coupons columns.
# BalaRaju - added code for departments

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two new columns -----------------------------------------
# New "HRA" column goes in before the existing "Spcl Allowance" column (old N, now shifts right).
$ws.Columns.Item(14).Insert()
# New "Food Coupons " column goes in before the existing "total_deducations" column
# (old U, which after the first insert sits at column 22).
$ws.Columns.Item(22).Insert()

# --- 2. Header row (row 1) --------------------------------------------------
$ws.Cells.Item(1, 14).Value = "HRA"
$ws.Cells.Item(1, 22).Value = "Food Coupons "

# --- 3. Data row (row 2) : refreshed figures + the two new columns ---------
$ws.Cells.Item(2, 8).Value  = 120000.0    # H2  GROSS
$ws.Cells.Item(2, 9).Value  = 10000.0     # I2  Per Month
$ws.Cells.Item(2, 10).Value = 9011.67     # J2  Actual Per Month
$ws.Cells.Item(2, 11).Value = 30.0        # K2  Actual Days
$ws.Cells.Item(2, 12).Value = 20.0        # L2  Working Days
$ws.Cells.Item(2, 13).Value = 2666.67     # M2  BASIC
$ws.Cells.Item(2, 14).Value = 800.0       # N2  HRA (new)
$ws.Cells.Item(2, 15).Value = 4545.0      # O2  Spcl Allowance
$ws.Cells.Item(2, 16).Value = 1000.0      # P2  Arrears
$ws.Cells.Item(2, 17).Value = 9011.67     # Q2  Gross Pay
$ws.Cells.Item(2, 18).Value = 320.0       # R2  PF
$ws.Cells.Item(2, 19).Value = 157.7       # S2  ESIC
$ws.Cells.Item(2, 20).Value = 100.0       # T2  PT
$ws.Cells.Item(2, 21).Value = 100.0       # U2  TDS
$ws.Cells.Item(2, 22).Value = 0           # V2  Food Coupons (new)
$ws.Cells.Item(2, 23).Value = 677.7       # W2  total_deducations
$ws.Cells.Item(2, 24).Value = 8333.97     # X2  NetPay

# --- 4. Column widths for the columns whose width actually changed --------
$ws.Columns.Item(8).ColumnWidth  = 10.142857142857142   # H
$ws.Columns.Item(9).ColumnWidth  = 10.142857142857142   # I
$ws.Columns.Item(13).ColumnWidth = 9.0                  # M
$ws.Columns.Item(14).ColumnWidth = 6.714285714285714    # N (new)
$ws.Columns.Item(16).ColumnWidth = 7.857142857142857    # P
$ws.Columns.Item(17).ColumnWidth = 9.0                  # Q
$ws.Columns.Item(18).ColumnWidth = 6.714285714285714    # R
$ws.Columns.Item(19).ColumnWidth = 6.714285714285714    # S
$ws.Columns.Item(20).ColumnWidth = 6.714285714285714    # T
$ws.Columns.Item(21).ColumnWidth = 6.714285714285714    # U
$ws.Columns.Item(22).ColumnWidth = 14.428571428571429   # V (new)
$ws.Columns.Item(24).ColumnWidth = 9.0                  # X

Write-Output "edit applied"
